$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-09-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-23 Monday", 2) | Out-Null

# Update each arithmetic expression cell by position (row-major),
# since some old/new values repeat across cells and a blind
# document-wide Find/Replace would mis-target them.
$t = $d.Tables.Item(1)
$values = @(
    "60-46=",
    "26+13=",
    "66-42=",
    "74+5=",
    "21+4=",
    "23-7=",
    "1+46=",
    "71-22=",
    "85-1=",
    "76-23=",
    "26+60=",
    "80-68=",
    "64-42=",
    "89-86=",
    "94-43=",
    "48+50=",
    "48+4=",
    "50+24=",
    "45-42=",
    "79-35=",
    "98-6=",
    "83+14=",
    "11+8=",
    "37-12=",
    "79+4=",
    "84-56=",
    "43+32=",
    "77-10=",
    "73-46=",
    "88-79=",
    "54+38=",
    "62-5=",
    "67-17=",
    "98-90=",
    "76-21=",
    "34+40=",
    "74-21=",
    "36+30=",
    "82-55=",
    "96-96=",
    "45-39=",
    "32-30=",
    "5+10=",
    "14-8=",
    "62+26=",
    "17+60=",
    "32+16=",
    "41+1=",
    "68+27=",
    "41-33=",
    "31+11=",
    "31+25=",
    "36-12=",
    "27-15=",
    "9+51=",
    "42+36=",
    "88-61=",
    "73+21=",
    "16+22=",
    "16+74=",
    "52-49=",
    "25-21=",
    "57-9=",
    "32+35=",
    "92-46=",
    "80-33=",
    "15+54=",
    "61-13=",
    "11+63=",
    "88-13=",
    "2+0=",
    "35+5=",
    "47+45=",
    "4+3=",
    "12+71=",
    "96-73=",
    "60+7=",
    "1+0=",
    "90-14=",
    "62+6=",
    "35-13=",
    "2+48=",
    "94-10=",
    "41+21=",
    "58+30=",
    "44-41=",
    "18+30=",
    "26+60=",
    "97-58=",
    "17+50=",
    "34-25=",
    "71+13=",
    "97-18=",
    "15+44=",
    "69-2=",
    "58-44=",
    "32+20=",
    "58+26=",
    "78-8=",
    "95-52="
)

$cols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}

Write-Host "Done updating $($values.Count) cells and the date heading."